# Updated cryptos list on Thu May 30 07:43:16 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin on Sheet1, and swaps row 51 from EnergySwap to Arweave (new
# name/link/price/volume), matching the latest scrape of the source data.
#
# Column D stores plain-text price strings (dotted thousands separators,
# fixed decimal places, tiny numbers in plain decimal form, etc.) rather
# than real numbers. Assigning a numeric-looking string straight to
# Range.Value lets Excel silently coerce it into a genuine number -
# dropping meaningful trailing zeros (e.g. "0.520" -> 0.52) or reformatting
# it (e.g. "0.0000264" -> 2.64E-05). A leading apostrophe is Excel's
# standard "force text" marker, so prefixing the literal with one keeps
# these values verbatim text, matching the source workbook's inline-string
# cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.916.08"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "'3.746.43"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'593.01"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").Value = "'165.38"
$ws.Range("E6").Value = "  -3.51%  "
$ws.Range("D7").Value = "'3.746.12"
$ws.Range("E7").Value = "  -2.27%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.520"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("D11").Value = "'6.43"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").Value = "'0.0000264"
$ws.Range("E13").Value = "  -6.79%  "
$ws.Range("D14").Value = "'35.91"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "'4.368.78"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").Value = "'3.742.92"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").Value = "'67.784.00"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "'18.34"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -5.59%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'10.59"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").Value = "'464.53"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'0.702"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").Value = "'82.91"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'0.0000136"
$ws.Range("E25").Value = "  -13.88%  "
$ws.Range("D26").Value = "'2.19"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("D27").Value = "'11.96"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "'3.887.19"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").Value = "'7.37"
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("D33").Value = "'29.88"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("D35").Value = "'9.03"
$ws.Range("E35").Value = "  -3.73%  "
$ws.Range("D36").Value = "'3.690.81"
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").Value = "'3.46"
$ws.Range("E38").Value = "  -11.12%  "
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'0.996"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").Value = "'5.75"
$ws.Range("E41").Value = "  -3.52%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -3.91%  "
$ws.Range("D45").Value = "'8.54"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").Value = "'1.91"
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("D47").Value = "'45.24"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("D48").Value = "'394.87"
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("D49").Value = "'144.70"
$ws.Range("D50").Value = "'0.0346"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").Value = "'38.32"
$ws.Range("E51").Value = "  +0.27%  "
